# Apply updated odds values to the FlashScore weekly games sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.1
$ws.Range("S2").Value = 1.5
$ws.Range("T2").Value = 2.37

# Row 3
$ws.Range("G3").Value = 2.1
$ws.Range("S3").Value = 1.54

# Row 4
$ws.Range("I4").Value = 2.3
$ws.Range("S4").Value = 1.58

# Row 5
$ws.Range("G5").Value = 1.83
$ws.Range("S5").Value = 1.47

# Row 6
$ws.Range("M6").Value = 1.11
$ws.Range("N6").Value = 6.5

# Row 7
$ws.Range("Q7").Value = 1.93
$ws.Range("R7").Value = 1.93

# Row 9
$ws.Range("M9").Value = 1.04
$ws.Range("N9").Value = 13
$ws.Range("O9").Value = 1.25
$ws.Range("P9").Value = 3.75
$ws.Range("Q9").Value = 1.83
$ws.Range("R9").Value = 2.03

# Row 10
$ws.Range("H10").Value = 3.6
$ws.Range("K10").Value = 2.3
$ws.Range("L10").Value = 3.5
$ws.Range("O10").Value = 1.2
$ws.Range("P10").Value = 4.33
$ws.Range("Q10").Value = 1.67
$ws.Range("R10").Value = 2.15
$ws.Range("U10").Value = 1.57
$ws.Range("V10").Value = 2.25
$ws.Range("W10").Value = 10
$ws.Range("AE10").Value = 12
$ws.Range("AG10").Value = 12
$ws.Range("AL10").Value = 26
$ws.Range("AM10").Value = 126
$ws.Range("AN10").Value = 4.5
$ws.Range("AP10").Value = 19
$ws.Range("AS10").Value = 101
$ws.Range("AV10").Value = 41
$ws.Range("AW10").Value = 5.5
$ws.Range("AY10").Value = 21
$ws.Range("BB10").Value = 126
